# Update the BOM sheets ("All", "Top", "Bottom") for the beta build:
#  - revision date A2: "2022 February 25" -> "2022 June 16"
#  - column header D4: "Mfgr Part #" -> "Name"
#  - PCB part number D9 (when present): "JS220_fp_bp_pcb_revB" -> "JS220_fp_binding_post_pcb_revC"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("A2").Value() -eq "2022 February 25") {
        $ws.Range("A2").Value = "2022 June 16"
    }

    if ($ws.Range("D4").Value() -eq "Mfgr Part #") {
        $ws.Range("D4").Value = "Name"
    }

    if ($ws.Range("D9").Value() -eq "JS220_fp_bp_pcb_revB") {
        $ws.Range("D9").Value = "JS220_fp_binding_post_pcb_revC"
    }
}
